# Slide 12 ("Records") - grammar listing for `variable`:
#   before:  ... paramId) { indexExpr | fieldExpr } ...
#   after:   ... paramId ) { indexExpr | fieldExpr } ...
# i.e. a single space is inserted right after "paramId" and before the
# closing paren, which (per the authored diff) lands as a brand-new run
# (" ) ") immediately followed by the remainder of the old run ("{ ").

$p = $ppt.ActivePresentation

# Locate the shape that holds the grammar text, searching every slide so
# this keeps working even if the slide order/index ever shifts.
$targetShape = $null
foreach ($sl in $p.Slides) {
    foreach ($sh in $sl.Shapes) {
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            if ($sh.TextFrame.TextRange.Text.IndexOf("paramId) { ") -ge 0) {
                $targetShape = $sh
                break
            }
        }
    }
    if ($targetShape -ne $null) { break }
}
if ($targetShape -eq $null) {
    throw "Could not find the shape containing the 'variable' production"
}

$tr = $targetShape.TextFrame.TextRange

# Grab the run that currently reads ") { " (right after the "paramId" run).
$run = $tr.Find(") { ", 0)
if ($run -eq $null) {
    throw "Could not find the ') { ' run to edit"
}

# Insert the missing space at the front: ") { " -> " ) { ".
$run.Text = " ) { "

# Split that text into two runs - " ) " and "{ " - matching the target
# markup, by nudging a (no-op) character formatting property on just the
# leading " ) " slice; touching a run like this is what makes the host
# materialize it as its own <a:r> while leaving the remaining "{ " text
# as a separate run with its original formatting intact.
$newRun = $tr.Characters($run.Start, 3)
$newRun.Font.Size = $newRun.Font.Size

Write-Host "Updated text:" $tr.Text
